$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Step 1: duplicate the existing "2022-Q2" sheet so we end up with two
# sheets: the original (which we will turn into "2022-Q3" with fresh
# data) and an untouched copy (which keeps the "2022-Q2" name/data).
# Worksheet.Copy(Before, After) with After = itself drops the copy
# right after the source, which is exactly where "2022-Q2" needs to
# land once the original slot is renamed to "2022-Q3".
# ------------------------------------------------------------------
$origQ2 = $wb.Worksheets.Item("2022-Q2")
$origQ2.Copy($null, $origQ2)

# the new copy gets auto-named "2022-Q2 (2)" since the name is taken
$copyQ2 = $wb.Worksheets.Item("2022-Q2 (2)")
$copyQ2.Name = "2022-Q2-staging"
$origQ2.Name = "2022-Q3"
$copyQ2.Name = "2022-Q2"

# ------------------------------------------------------------------
# Step 2: update the "总计" (summary) sheet - insert the new 2022-Q3
# row at the top of the data and push the existing rows down one.
# ------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Pick up row 4's style (border/alignment) for the brand-new row 5
# before overwriting any values.
$summary.Range("A4").Copy($summary.Range("A5"))

# Shift rows downward, bottom-up so we never read an already-clobbered
# source cell.
$summary.Range("B5").Value = $summary.Range("B4").Value2
$summary.Range("C5").Value = $summary.Range("C4").Value2
$summary.Range("D5").Value = $summary.Range("D4").Value2
$summary.Range("A5").Value = 3

$summary.Range("B4").Value = $summary.Range("B3").Value2
$summary.Range("C4").Value = $summary.Range("C3").Value2
$summary.Range("D4").Value = $summary.Range("D3").Value2

$summary.Range("B3").Value = $summary.Range("B2").Value2
$summary.Range("C3").Value = $summary.Range("C2").Value2
$summary.Range("D3").Value = $summary.Range("D2").Value2

$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 1
$summary.Range("D2").Value = 3.01

# ------------------------------------------------------------------
# Step 3: fill in the real 2022-Q3 figures on the new "2022-Q3" sheet.
# ------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")

$q3.Range("C2").Value = "易方达亚洲精选股票（QDII）"

# D2/E2/F2/G2 hold numeric-looking text (not real numbers) in the
# source file, so force text formatting, write the value, then drop
# the formatting again so no stray style index is left behind.
$q3.Range("D2").NumberFormat = "@"
$q3.Range("D2").Value = "46.17"
$q3.Range("D2").ClearFormats()

$q3.Range("E2").NumberFormat = "@"
$q3.Range("E2").Value = "94.52"
$q3.Range("E2").ClearFormats()

$q3.Range("F2").NumberFormat = "@"
$q3.Range("F2").Value = "6.51"
$q3.Range("F2").ClearFormats()

$q3.Range("G2").NumberFormat = "@"
$q3.Range("G2").Value = "3.0057"
$q3.Range("G2").ClearFormats()

$q3.Range("H2").Value = 4

# ------------------------------------------------------------------
# Step 4: the sheet-copy above shifted the "active" tab; restore it to
# "2021-Q2" (the last sheet), which is where it was selected before.
# ------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q2").Activate()
